# Append one new row (row 9) to Sheet1, matching the row shape already
# used by rows 1-8 (every cell is plain text, even the numeric-looking
# ones such as the "quantity" column C).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 9

# Column A ("ملاحظات" / notes) is blank for this entry, same as most
# other rows.
$ws.Cells.Item($row, 1).Value = ""

$ws.Cells.Item($row, 2).Value = "احمد"

# Column C ("الكمية" / quantity) holds a numeric-looking value but must
# stay stored as text (like the rest of the sheet, which is flagged with
# numberStoredAsText). Force text formatting before assigning the value
# so it isn't auto-converted to a number, then drop the now-unneeded
# number-format override so no stray style is left behind.
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "23"
$ws.Cells.Item($row, 3).ClearFormats()

$ws.Cells.Item($row, 4).Value = "الصمود"
$ws.Cells.Item($row, 5).Value = "الرحلة 1"
$ws.Cells.Item($row, 6).Value = "C1"
$ws.Cells.Item($row, 7).Value = "UNICEF"
$ws.Cells.Item($row, 8).Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٦:٥١:٠٥ م"
